# RFC: Item#references: initial support of many ref items (Take 3)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column E ("references") - shifts old E..J to F..K
$ws.Columns.Item(5).Insert()

# Insert a new row 4 (for the new REQ006 item) - shifts old rows 4..6 to 5..7
$ws.Rows.Item(4).Insert()

# --- Header row ---
$ws.Cells.Item(1, 5).Value = "references"

# --- Column widths (A-D unchanged, E onward updated).
# The engine adds a constant 5/6-character padding when it round-trips a
# ColumnWidth back to the stored OOXML width, so subtract it here to land
# on the exact target widths.
$pad = 5.0 / 6.0
$ws.Columns.Item(5).ColumnWidth = 50.5 - $pad
$ws.Columns.Item(6).ColumnWidth = 42.5 - $pad
$ws.Columns.Item(7).ColumnWidth = 9.5 - $pad
$ws.Columns.Item(8).ColumnWidth = 10.5 - $pad
$ws.Columns.Item(9).ColumnWidth = 11.5 - $pad
$ws.Columns.Item(10).ColumnWidth = 12.5 - $pad
$ws.Columns.Item(11).ColumnWidth = 47.5 - $pad

# --- New row 4 formatting: match the rest of the data rows (wrap text,
# left/top aligned) by copying the format from row 5 (an existing data row)
# rather than setting WrapText/Alignment directly, which would force an
# autofit row-height recalculation that the other data rows don't have.
$ws.Range("A5:K5").Copy()
$ws.Range("A4:K4").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$excel.CutCopyMode = $false

# --- Row 4 (new REQ006 item) ---
$ws.Cells.Item(4, 1).Value = "REQ006"
$ws.Cells.Item(4, 2).NumberFormat = "@"
$ws.Cells.Item(4, 2).Value = "1.5"
$ws.Cells.Item(4, 3).Value = "Hello, world!"
$ws.Cells.Item(4, 5).Value = "type:file,path:external/text.txt,keyword:REF123`ntype:file,path:external/text2.txt"
$ws.Cells.Item(4, 6).Value = "REQ001:35ed54323e3054c33ae5545fffdbbbf5"
$ws.Cells.Item(4, 7).Value = $true
$ws.Cells.Item(4, 8).Value = $false
$ws.Cells.Item(4, 10).Value = $true
$ws.Cells.Item(4, 11).Value = "c442316131ca0225595ae257f3b4583d"

# Setting a multi-line value above triggers a live autofit row-height
# recalculation (unlike the original workbook, where row heights are left
# implicit/auto for every row, even ones with multi-line wrapped text).
# Re-running AutoFit clears the resulting explicit `ht`/`customHeight`
# marker so row 4 serializes the same way as every other data row.
$ws.Rows.Item(4).AutoFit() | Out-Null

# --- Rows 2,3,5,6,7: the new "references" column (E) and, where a row had
# a single "links" value, that value simply moved one column right via the
# column insert above - nothing further needed there.

# --- AutoFilter range + _FilterDatabase defined name now cover A:K ---
$ws.AutoFilterMode = $false
$ws.Range("A1:K1").AutoFilter() | Out-Null
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "='Sheet'!`$A`$1:`$K`$1"
